$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.930.43'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.89%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.262.78'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.28%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('B5').Value = 'XRP'
$ws.Range('C5').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.654'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.45%  '

$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '233.40'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.17%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '63.86'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.03%  '

$ws.Range('E8').Value = '  -0.14%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.449'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.20%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0977'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.91%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '58.31'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.87%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '26.55'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.80%  '

$ws.Range('E13').Value = '  +1.47%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.600.14'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.35%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.63'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.08%  '

$ws.Range('E16').Value = '  +4.42%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.843'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.77%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.260.17'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.66%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '43.838.47'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.76%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0979'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.10%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.79'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.75%  '

$ws.Range('E22').Value = '  +1.54%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '249.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.30%  '

$ws.Range('E24').Value = '  -0.01%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.46'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.85%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.57'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +27.68%  '

$ws.Range('E27').Value = '  -3.47%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.91'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.83%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '174.05'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.46%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '21.90'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.80%  '

$ws.Range('E31').Value = '  +0.08%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.43'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.03%  '

$ws.Range('E33').Value = '  +3.45%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.97'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.59%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0686'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.27%  '

$ws.Range('E36').Value = '  -1.64%  '

$ws.Range('E37').Value = '  -2.12%  '

$ws.Range('E38').Value = '  -5.64%  '

$ws.Range('E39').Value = '  -1.17%  '

$ws.Range('E40').Value = '  +3.27%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.02%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.79'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.02%  '

$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.41'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.07%  '

$ws.Range('B44').Value = 'FTXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.50'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.80%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '98.85'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.68%  '

$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0952'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.81%  '

$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.20'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.45%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000210'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.48%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.455.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.37%  '

$ws.Range('B50').Value = 'Celestia'
$ws.Range('C50').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.46%  '

$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.32'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.44%  '
